# Add "kamal" as a new employee row (row 3) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New employee data
$ws.Range("A3").Value = 102
$ws.Range("B3").Value = "kamal"
$ws.Range("C3").Value = "kamal@gmail.com"

# Add the mailto hyperlink on the email cell (same as existing salam row),
# then re-apply the Hyperlink cell style so it matches C2 exactly.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:kamal@gmail.com")
$ws.Range("C3").Style = "Hyperlink"

# Widen column C to fit the new email text (matches Excel's auto column sizing)
$ws.Columns.Item(3).ColumnWidth = 21.6

# The active cell/selection moves on to C4 after entering the new row
$ws.Range("C4").Select()
